$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly records were collected for "Repollo / Crespo record" at
# "Macroferia Regional de Talca". Insert two new rows above the existing
# row 262 (pushing the rest of the table, currently rows 262:281, down to
# 264:283) and populate them with the new observations.
$ws.Rows("262:263").Insert()

# New row 262: Primera calidad, fecha 44706
$ws.Range("A262").Value = 5
$ws.Range("B262").Value = "Macroferia Regional de Talca"
$ws.Range("C262").Value = "Maule"
$ws.Range("D262").Value = 44706
$ws.Range("E262").Value = 7
$ws.Range("F262").Value = 100112006
$ws.Range("G262").Value = "Repollo"
$ws.Range("H262").Value = "Crespo record"
$ws.Range("I262").Value = "Primera"
$ws.Range("J262").Value = 2000
$ws.Range("K262").Value = 1100
$ws.Range("L262").Value = 1100
$ws.Range("M262").Value = 1100
$ws.Range("N262").Value = '$/unidad'
$ws.Range("O262").Value = "Región del Maule"
$ws.Range("P262").Value = 1100
$ws.Range("Q262").Value = 1
$ws.Range("R262").Value = "Hortaliza"

# New row 263: Segunda calidad, fecha 44706
$ws.Range("A263").Value = 5
$ws.Range("B263").Value = "Macroferia Regional de Talca"
$ws.Range("C263").Value = "Maule"
$ws.Range("D263").Value = 44706
$ws.Range("E263").Value = 7
$ws.Range("F263").Value = 100112006
$ws.Range("G263").Value = "Repollo"
$ws.Range("H263").Value = "Crespo record"
$ws.Range("I263").Value = "Segunda"
$ws.Range("J263").Value = 2000
$ws.Range("K263").Value = 900
$ws.Range("L263").Value = 900
$ws.Range("M263").Value = 900
$ws.Range("N263").Value = '$/unidad'
$ws.Range("O263").Value = "Región del Maule"
$ws.Range("P263").Value = 900
$ws.Range("Q263").Value = 1
$ws.Range("R263").Value = "Hortaliza"
